$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying simulation output changed: two out-of-sequence "Kairi94..."
# observation rows (serial date 34213) were dropped - one from the
# Nitrogen100 block (row 234) and one from the Nitrogen150 block (row 246,
# which becomes row 245 once the first row above it has been removed).
# Deleting the rows shifts everything below them up by one, which reproduces
# the rest of the row renumbering seen in the diff.
$ws.Rows("234:234").Delete()
$ws.Rows("245:245").Delete()

# Refresh the AutoFilter range now that the sheet is two rows shorter.
$ws.AutoFilterMode = $false
$ws.Range("A1:A341").AutoFilter() | Out-Null

# Keep the hidden _FilterDatabase defined name in sync with the new range.
$wb.Names.Item("Observed!_FilterDatabase").RefersTo = "=Observed!`$A`$1:`$A`$341"

# Match the author's final view state: scrolled down near the bottom of the
# (now shorter) sheet, with D249 selected.
$excel.ActiveWindow.ScrollRow = 216
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D249").Select() | Out-Null
